# Capacidade.xlsx - "Add files via upload"
#
# The sheet "Placas" has a license-plate entry in A2 and a schedule
# (Escala) chosen from a dropdown (data validation list) in B2. The
# commit swaps both for a new plate/schedule pair; D2 (VLOOKUP of B2
# against Quantidade!A:B) and the Quantidade!C4/C5 duplicate counters
# recompute automatically. Finally the sheet's last active cell moves
# to B5.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Placas")
$ws.Activate()

# New plate value (was RNS0H99).
$ws.Range("A2").Value = "BEP2H98"

# Pick a different item from the Escala dropdown (was "Diarista 8 - 14").
$ws.Range("B2").Value = "Diarista 6 - 12"

# Last selection left on the sheet after the edit.
$ws.Range("B5").Select()
